$d = $word.ActiveDocument

# 1. Fix "START TRANSCATION;" typo -> "START TRANSACTION;"
$d.Content.Find.Execute("START TRANSCATION;", $true, $false, $false, $false, $false,
                         $true, 1, $false, "START TRANSACTION;", 2)

# 2. Remove comma: "SELECT keyword, is classified" -> "SELECT keyword is classified"
$d.Content.Find.Execute("the SELECT keyword, is classified as an element of DML.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "the SELECT keyword is classified as an element of DML.", 2)

# 3. Subject/verb agreement: "missing data are tolerated by MySQL" -> "missing data is tolerated by MySQL"
$d.Content.Find.Execute("missing data are tolerated by MySQL", $true, $false, $false, $false, $false,
                         $true, 1, $false, "missing data is tolerated by MySQL", 2)

# 4. Remove "14. " numbering prefix from the "Open Notepad" paragraph
$d.Content.Find.Execute("14. Open Notepad. Repeat Exercise (6), typing the INSERT statements in Notepad. Save the file as ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Open Notepad. Repeat Exercise (6), typing the INSERT statements in Notepad. Save the file as ", 2)

# 5. Renumber "15. Execute the" -> "14. Execute the"
$d.Content.Find.Execute("15. Execute the", $true, $false, $false, $false, $false,
                         $true, 1, $false, "14. Execute the", 2)

# 6. Renumber "16. Issue a SAVEPOINT command. Name the " -> "15. Issue a SAVEPOINT command. Name the "
$d.Content.Find.Execute("16. Issue a SAVEPOINT command. Name the ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "15. Issue a SAVEPOINT command. Name the ", 2)

# 7. Renumber "17. Write statements to make the following UPDATEs to " -> "16. Write statements to make the following UPDATEs to "
$d.Content.Find.Execute("17. Write statements to make the following UPDATEs to ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "16. Write statements to make the following UPDATEs to ", 2)

# 8. "(1)" -> "(a)"
$d.Content.Find.Execute("(1)", $true, $false, $false, $false, $false,
                         $true, 1, $false, "(a)", 2)

# 9. "(2)" -> "(b)"
$d.Content.Find.Execute("(2)", $true, $false, $false, $false, $false,
                         $true, 1, $false, "(b)", 2)

# 10. "(3)" -> "(c)"
$d.Content.Find.Execute("(3)", $true, $false, $false, $false, $false,
                         $true, 1, $false, "(c)", 2)

# 11. Renumber "18. Write a SELECT statement..." -> "17. Write a SELECT statement..."
$d.Content.Find.Execute("18. Write a SELECT statement to display the Shipper and each modified field to verify the changes.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "17. Write a SELECT statement to display the Shipper and each modified field to verify the changes.", 2)

# 12. Renumber "19. Rollback the changes made to your " -> "18. Rollback the changes made to your "
$d.Content.Find.Execute("19. Rollback the changes made to your ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "18. Rollback the changes made to your ", 2)

# 13. Renumber "20. Commit the changes from this transaction." -> "19. Commit the changes from this transaction."
$d.Content.Find.Execute("20. Commit the changes from this transaction.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "19. Commit the changes from this transaction.", 2)

# 14. Renumber "21. Write DELETE statements to remove the records " -> "20. Write DELETE statements to remove the records "
$d.Content.Find.Execute("21. Write DELETE statements to remove the records ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "20. Write DELETE statements to remove the records ", 2)

# 15. Renumber "22. " -> "21. " (before "Cost_table should now be back")
$d.Content.Find.Execute("22. Cost_table should now be back", $true, $false, $false, $false, $false,
                         $true, 1, $false, "21. Cost_table should now be back", 2)

Write-Host "Done"
